$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 26-27, pushing existing rows 26.. down to 28..
$ws.Rows("26:27").Insert()

# --- Fill in new row 26 ---
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 44742
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100106
$ws.Cells.Item(26, 8).Value = "Oleaginosos"
$ws.Cells.Item(26, 9).Value = 100106002
$ws.Cells.Item(26, 10).Value = "Palta"
$ws.Cells.Item(26, 11).Value = "Fuerte"
$ws.Cells.Item(26, 12).Value = "Tercera"
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 42000
$ws.Cells.Item(26, 15).Value = 43000
$ws.Cells.Item(26, 16).Value = 42500
$ws.Cells.Item(26, 17).Value = "$/caja 25 kilos"
$ws.Cells.Item(26, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(26, 19).Value = 1700
$ws.Cells.Item(26, 20).Value = 25

# --- Fill in new row 27 ---
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44742
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100106
$ws.Cells.Item(27, 8).Value = "Oleaginosos"
$ws.Cells.Item(27, 9).Value = 100106002
$ws.Cells.Item(27, 10).Value = "Palta"
$ws.Cells.Item(27, 11).Value = "Hass"
$ws.Cells.Item(27, 12).Value = "Segunda"
$ws.Cells.Item(27, 13).Value = 200
$ws.Cells.Item(27, 14).Value = 14000
$ws.Cells.Item(27, 15).Value = 15000
$ws.Cells.Item(27, 16).Value = 14500
$ws.Cells.Item(27, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(27, 18).Value = "Perú"
$ws.Cells.Item(27, 19).Value = 1450
$ws.Cells.Item(27, 20).Value = 10
